$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "ordin" + (old _GoBack bookmark) + "ale"  ->  single run "ordinale"
#    A same-text Find & Replace collapses the two runs that straddle
#    the old bookmark into a single run and drops the now-orphaned
#    _GoBack bookmark that used to sit between them.
# ------------------------------------------------------------------
$rOrdinale = $d.Content
$rOrdinale.Find.Execute("ordinale", $false, $false, $false, $false, $false, `
    $true, 1, $false, "ordinale", 2)

# ------------------------------------------------------------------
# 2) "Abovementioned the tree classifier can also be used"
#    -> "Above" + " " + "mentioned the tree classifier can also be used"
#    split across three runs (grammar fix: "Abovementioned" -> "Above
#    mentioned"), while leaving the following runs untouched.
# ------------------------------------------------------------------
$rAbove = $d.Content
$rAbove.Find.Execute("Abovementioned the tree classifier can also be used")
$aboveStart = $rAbove.Start
$aboveEnd = $rAbove.End

# Insert the missing space right after "Above" (5 characters in).
$insertionPoint = $d.Range($aboveStart + 5, $aboveStart + 5)
$insertionPoint.InsertAfter(" ")

# Re-assert (identical) formatting on each of the three desired runs so
# the engine keeps them as separate <w:r> elements instead of folding
# them back into one contiguous run when the package is saved.
$runAbove = $d.Range($aboveStart, $aboveStart + 5)
$runAbove.Font.Bold = 1
$runAbove.Font.Bold = 0

$runSpace = $d.Range($aboveStart + 5, $aboveStart + 6)
$runSpace.Font.Bold = 1
$runSpace.Font.Bold = 0

$runMentioned = $d.Range($aboveStart + 6, $aboveEnd + 1)
$runMentioned.Font.Bold = 1
$runMentioned.Font.Bold = 0

# Re-assert the formatting of the still-untouched tail run as well so
# that it keeps its own boundary instead of being folded together with
# the "mentioned ..." run above when the package is normalised/saved.
$rTail = $d.Content
$rTail.Find.Execute(" for out-of-sample predictions and with reference to its accuracy 78.89% it can be stated out that it is around 14% better than the defined baseline.")
$runTail = $d.Range($rTail.Start, $rTail.End)
$runTail.Font.Bold = 1
$runTail.Font.Bold = 0

# ------------------------------------------------------------------
# 3) Move the _GoBack bookmark so it wraps the whole paragraph that
#    starts with "The predictive analysis ..." and ends right after
#    "... defined baseline. " (i.e. the paragraph containing the text
#    edited in step 2).
# ------------------------------------------------------------------
$rPara = $d.Content
$rPara.Find.Execute("The predictive analysis is in contrast")
$paraStart = $rPara.Start

$rParaEnd = $d.Content
$rParaEnd.Find.Execute("around 14% better than the defined baseline.")
$paraContentEnd = $rParaEnd.End

# Extend to include the trailing space run, but stop before the
# paragraph mark itself.
$paraRange = $d.Range($paraStart, $paraContentEnd + 1)
$d.Bookmarks.Add("_GoBack", $paraRange)
